# Generate Report for Handback
# Update the "generated date" / handoff / handback timestamp cells on the
# Overview, zh-cn and de-de sheets to reflect a newer report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on the Overview sheet.
$wsOverview.Range("G2").Value = "2016-08-31 15:24:29"

# "Correspond Handoff Datetime" / "Correspond Handback DateTime" on zh-cn.
$wsZhCn.Range("H2").Value = "2016-08-31 15:24:24"
$wsZhCn.Range("K2").Value = "2016-08-31 15:24:42"

# de-de shares the same "Latest HO Xliff Generate Date" value as Overview,
# plus its own "Correspond Handback DateTime".
$wsDeDe.Range("H2").Value = "2016-08-31 15:24:29"
$wsDeDe.Range("K2").Value = "2016-08-31 15:24:49"
